$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("existing_stock")

$ws.Range('E14').Value = 0.066
$ws.Range('H14').Value = 60.500000000000014
$ws.Range('E15').Value = 0.086
$ws.Range('H15').Value = 60.50000000000001
$ws.Range('C52').Value = 'e_w207993342-220'
$ws.Range('C82').Value = 'e_w1284913429-220'
$ws.Range('E82').Value = 0.05
$ws.Range('G82').Value = 3267.0000000000005
$ws.Range('H82').Value = 78.65
$ws.Range('I82').Value = 3.1500000000000004
$ws.Range('C83').Value = 'e_w391576135-220'
$ws.Range('E83').Value = 0.061
$ws.Range('G83').Value = 2783.0
$ws.Range('H83').Value = 66.55000000000001
$ws.Range('I83').Value = 2.8875
$ws.Range('E95').Value = 0.001
$ws.Range('C96').Value = 'e_w234983117-220'
$ws.Range('E96').Value = 0.001
$ws.Range('C97').Value = 'e_w27435934-220'
$ws.Range('C98').Value = 'e_w89977424-220'
$ws.Range('E98').Value = 0.0011
$ws.Range('E99').Value = 0.0012
$ws.Range('C101').Value = 'e_w234983117-220'
$ws.Range('E101').Value = 0.0023
$ws.Range('G101').Value = 1336.5
$ws.Range('H101').Value = 21.450000000000006
$ws.Range('C102').Value = 'e_w97941869-220'
$ws.Range('E102').Value = 0.0015
$ws.Range('C103').Value = 'e_CH60-225'
$ws.Range('E103').Value = 0.0034
$ws.Range('G103').Value = 1336.5000000000002
$ws.Range('E104').Value = 0.003
$ws.Range('C105').Value = 'e_w89977424-220'
$ws.Range('E105').Value = 0.0015
$ws.Range('G105').Value = 1336.5
$ws.Range('E106').Value = 0.001
$ws.Range('E107').Value = 0.0018
$ws.Range('G107').Value = 1336.5000000000002
$ws.Range('C108').Value = 'e_w281809991-220'
$ws.Range('E108').Value = 0.0012
$ws.Range('C109').Value = 'e_w1105061707-220'
$ws.Range('E109').Value = 0.0016
$ws.Range('H109').Value = 21.450000000000003
$ws.Range('C110').Value = 'e_w127004407-380'
$ws.Range('E110').Value = 0.0017
$ws.Range('G110').Value = 1336.4999999999998
$ws.Range('C111').Value = 'e_CH57-220'
$ws.Range('E111').Value = 0.0058000000000000005
$ws.Range('C112').Value = 'e_CH17-380'
$ws.Range('E112').Value = 0.0012
$ws.Range('C113').Value = 'e_w1105061707-220'
$ws.Range('C114').Value = 'e_w190819048-220'
$ws.Range('E114').Value = 0.001
$ws.Range('C115').Value = 'e_CH11-220'
$ws.Range('E115').Value = 0.0070999999999999995
$ws.Range('G115').Value = 1336.5000000000002
$ws.Range('C116').Value = 'e_w109037817-220'
$ws.Range('E116').Value = 0.0021000000000000003
$ws.Range('C117').Value = 'e_w281809991-220'
$ws.Range('E117').Value = 0.0013
$ws.Range('C118').Value = 'e_CH31-220'
$ws.Range('E118').Value = 0.0012
$ws.Range('G118').Value = 1336.5
$ws.Range('C119').Value = 'e_CH60-225'
$ws.Range('E119').Value = 0.0011
$ws.Range('C121').Value = 'e_w97941869-220'
$ws.Range('E121').Value = 0.0015
$ws.Range('C122').Value = 'e_w89977424-220'
$ws.Range('E122').Value = 0.0011
$ws.Range('C123').Value = 'e_w55695765-220'
$ws.Range('E123').Value = 0.001
$ws.Range('C124').Value = 'e_w165254212-220'
$ws.Range('E124').Value = 0.0058
$ws.Range('C125').Value = 'e_w1105061707-220'
$ws.Range('E125').Value = 0.0045
$ws.Range('C126').Value = 'e_CH11-220'
$ws.Range('E126').Value = 0.0012
$ws.Range('C127').Value = 'e_CH17-380'
$ws.Range('E127').Value = 0.0021000000000000003
$ws.Range('C128').Value = 'e_r5378910-220'
$ws.Range('E128').Value = 0.0013
$ws.Range('C129').Value = 'e_w281809991-220'
$ws.Range('E129').Value = 0.005
$ws.Range('C130').Value = 'e_w127004407-380'
$ws.Range('E130').Value = 0.004900000000000001
$ws.Range('C132').Value = 'e_CH31-220'
$ws.Range('C133').Value = 'e_w97941869-220'
$ws.Range('E133').Value = 0.0022
$ws.Range('E135').Value = 0.21471510601685545
$ws.Range('E136').Value = 0.18308354646436523
$ws.Range('E137').Value = 0.211523174241075
$ws.Range('E138').Value = 0.16277902359433066
$ws.Range('E139').Value = 0.1982862967966156
$ws.Range('E140').Value = 0.1586174511333161
$ws.Range('E141').Value = 0.15521278721895346
$ws.Range('E142').Value = 0.16554834018408843
$ws.Range('E143').Value = 0.15291072157643879
$ws.Range('E144').Value = 0.1727926292604506
$ws.Range('E145').Value = 0.1930981018275324
$ws.Range('E146').Value = 0.1663793251004252
$ws.Range('E147').Value = 0.19782569372870323
$ws.Range('E148').Value = 0.21077006448261207
$ws.Range('E149').Value = 0.21731537653220406
$ws.Range('E150').Value = 0.16699472878703805
$ws.Range('E151').Value = 0.13606784444360143
$ws.Range('E152').Value = 0.19614947844032105
$ws.Range('E153').Value = 0.1533817713118708
$ws.Range('E154').Value = 0.19328994063107527
$ws.Range('E155').Value = 0.1371001651339535
$ws.Range('E156').Value = 0.15409499379434963
$ws.Range('E157').Value = 0.20091315882928704
$ws.Range('E158').Value = 0.16152827258311295
$ws.Range('E159').Value = 0.1708220078874234
$ws.Range('C160').Value = 'e_w127004407-380'
$ws.Range('C162').Value = 'e_w127004407-380'
$ws.Range('C164').Value = 'e_CH31-220'
